# Add the 3 new books that were ranked to the bottom of the "Book Ranking" list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(72, 1).Value = "Bringin up Bebe"
$ws.Cells.Item(72, 2).Value = 4
$ws.Cells.Item(72, 3).Value = 432
$ws.Cells.Item(72, 4).Formula = "=B72/C72"
$ws.Cells.Item(72, 5).Value = "No"

$ws.Cells.Item(73, 1).Value = "Fiasco"
$ws.Cells.Item(73, 2).Value = 4
$ws.Cells.Item(73, 3).Value = 512
$ws.Cells.Item(73, 4).Formula = "=B73/C73"
$ws.Cells.Item(73, 5).Value = "No"

$ws.Cells.Item(74, 1).Value = "The Relativity of Wrong (Asimov)"
$ws.Cells.Item(74, 2).Value = 2
$ws.Cells.Item(74, 3).Value = 225
$ws.Cells.Item(74, 4).Formula = "=B74/C74"
$ws.Cells.Item(74, 5).Value = "No"

# Re-sort the whole table (rows 3-74) by the Result column (D), descending,
# exactly like the original list was sorted.
$rng = $ws.Range("A3:E74")
$key = $ws.Range("D3")
$rng.Sort($key, 2, $null, $null, 1, $null, 1, 1)

# Match the workbook's last on-screen selection after the edit.
$ws.Range("E74").Select()
